# Add "activity_aim" column (P) to the lesson plan sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Header for the new column P, mirrors the existing header style/shared string usage.
$ws.Cells.Item(1, 16).Value2 = "activity_aim"

# Data row: reuse the same placeholder text already used by the neighboring
# column O (row 2) so it maps onto the same shared string.
$oValue = $ws.Cells.Item(2, 15).Value2
$ws.Cells.Item(2, 16).Value2 = $oValue

# Update the view: scroll so column I is the left-most visible column, and
# select P2 (mirrors the saved selection/top-left cell in the workbook).
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("P2").Select()

$wb.Save()
